$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 6 with the values that used to be in row 8
$ws.Range("A6").Value = 7
$ws.Range("B6").Value = 9

# Delete rows 7 and 8 (the old row 6/7 data no longer needed)
$ws.Range("A7:B8").EntireRow.Delete()

# Update the active selection to match the target state
$ws.Range("C7").Select()
